# Edit 1: Insert new ToDo list item after "Download chart data and load to database / view"
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Download chart data and load to database / view")
if (-not $found) {
    throw "Could not find target paragraph 'Download chart data and load to database / view'"
}
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# Re-locate the paragraph we just split and grab the freshly inserted (empty) one after it.
$anchorPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -eq "Download chart data and load to database / view`r") {
        $anchorPara = $para
    }
}
if ($anchorPara -eq $null) {
    throw "Could not re-locate anchor paragraph after InsertParagraphAfter"
}
$newPara = $anchorPara.Next()

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Later: Code </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>thunk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ‘getting data’ actions to display to user that data is being downloaded…</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newParaXml)

# Edit 2 & 3: lastRenderedPageBreak moves from the "Stock Selection Form" run to the "NavBar" run
# (a side effect of the new paragraph above pushing the page break up by one item).

# 2a: add <w:lastRenderedPageBreak/> to the run holding "NavBar" (the list item right before
#     "Stock Selection Form", under the Features / User Experience numbering, numId=2).
$navRng = $d.Content
$navFound = $navRng.Find.Execute("NavBar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $navFound) {
    throw "Could not find 'NavBar' run"
}
$navPara = $navRng.Paragraphs(1)
$navParaRange = $navPara.Range
$navXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="002455F5" w:rsidRDefault="002455F5" w:rsidP="00B245EF"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>NavBar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> that </w:t></w:r><w:r w:rsidR="001026B3"><w:t>dynamically</w:t></w:r><w:r><w:t xml:space="preserve"> updates various data and information</w:t></w:r></w:p>'
$navParaRange.InsertXML($navXml)

# 2b: remove <w:lastRenderedPageBreak/> from the run holding "Stock Selection Form".
$stockRng = $d.Content
$stockFound = $stockRng.Find.Execute("Stock Selection Form", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $stockFound) {
    throw "Could not find 'Stock Selection Form' paragraph"
}
$stockPara = $stockRng.Paragraphs(1)
$stockParaRange = $stockPara.Range
$stockXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="002455F5" w:rsidRDefault="002455F5" w:rsidP="00B245EF"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Stock Selection Form</w:t></w:r></w:p>'
$stockParaRange.InsertXML($stockXml)

Write-Output "Edits applied successfully"
